$wb = $excel.ActiveWorkbook

# --- Sheet: VENTAS POR GRUPO ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("D20").Value = 274.75
$ws1.Range("D24").Value = "2 de 22"

# --- Sheet: VENTA MENSUAL ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F20").Value = 274.75
$ws2.Range("F24").Value = 4032.2

# --- Sheet: CUMPLIMIENTO MENSUAL ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D3").Value = 559.6799999999999
$ws3.Range("E3").Value = 8109.23
$ws3.Range("F3").Value = 0.06456174997779421
$ws3.Range("D19").Value = 4032.2
$ws3.Range("E19").Value = 50990.96386304604
$ws3.Range("F19").Value = 0.07328186379896733
